# Updates cryptos list values per the commit diff (prices + volume %).
# Numeric-looking price strings are forced to Text format first so Excel
# preserves them verbatim (matching trailing zeros / multi-dot values)
# instead of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.169.71'
$ws.Range("E2").Value = '  -0.52%  '

# Row 3
$ws.Range("D3").Value = '1.863.51'
$ws.Range("E3").Value = '  -0.41%  '

# Row 4
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.74'
$ws.Range("E5").Value = '  -1.15%  '

# Row 6
$ws.Range("E6").Value = '  -0.04%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4668'
$ws.Range("E7").Value = '  -0.95%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.13'
$ws.Range("E8").Value = '  +0.35%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2854'
$ws.Range("E9").Value = '  -1.50%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06477'
$ws.Range("E10").Value = '  -2.23%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.21'
$ws.Range("E11").Value = '  -2.15%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07736'
$ws.Range("E12").Value = '  -3.95%  '

# Row 13
$ws.Range("D13").Value = '1.859.03'
$ws.Range("E13").Value = '  -0.62%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.70'
$ws.Range("E14").Value = '  -3.76%  '

# Row 15
$ws.Range("E15").Value = '  -1.03%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.037'
$ws.Range("E16").Value = '  -2.02%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '268.53'
$ws.Range("E17").Value = '  -1.36%  '

# Row 18
$ws.Range("D18").Value = '30.153.08'
$ws.Range("E18").Value = '  -0.57%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.30'
$ws.Range("E19").Value = '  -6.02%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007594'
$ws.Range("E20").Value = '  -1.69%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.06%  '

# Row 22
$ws.Range("D22").Value = '2.081.63'
$ws.Range("E22").Value = '  -1.70%  '

# Row 23
$ws.Range("E23").Value = '  -0.05%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.137'
$ws.Range("E24").Value = '  -3.43%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.092'
$ws.Range("E25").Value = '  -2.17%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.311'
$ws.Range("E26").Value = '  -0.08%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.46'
$ws.Range("E27").Value = '  -1.43%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.58'
$ws.Range("E28").Value = '  -2.18%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.887'
$ws.Range("E29").Value = '  -3.66%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.361'
$ws.Range("E30").Value = '  -0.82%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09844'
$ws.Range("E31").Value = '  -1.38%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.449'
$ws.Range("E32").Value = '  -1.01%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.227'
$ws.Range("E33").Value = '  -3.33%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.005'
$ws.Range("E34").Value = '  -2.09%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.04674'
$ws.Range("E35").Value = '  -0.76%  '

# Row 36
$ws.Range("E36").Value = '  -1.49%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6874'
$ws.Range("E37").Value = '  -2.23%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.706'
$ws.Range("E38").Value = '  -0.28%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01829'
$ws.Range("E39").Value = '  -3.05%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.741'
$ws.Range("E40").Value = '  +3.41%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.279'
$ws.Range("E41").Value = '  -0.55%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.93'
$ws.Range("E42").Value = '  -2.80%  '

# Row 43
$ws.Range("E43").Value = '  +0.00%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.891'
$ws.Range("E44").Value = '  -3.77%  '

# Row 45
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8334'
$ws.Range("E45").Value = '  -1.10%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.93'
$ws.Range("E46").Value = '  -1.38%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4051'
$ws.Range("E47").Value = '  -3.01%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '934.59'
$ws.Range("E48").Value = '  -0.15%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.053'
$ws.Range("E49").Value = '  -2.76%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.948'
$ws.Range("E50").Value = '  -2.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.92'
$ws.Range("E51").Value = '  -1.87%  '
